# Apply cryptos list update (prices and 1h volume %) per commit:
# "Updated cryptos list on Thu May  2 16:30:26 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.130.69'
$ws.Range("E2").Value = '  +3.30%  '
$ws.Range("D3").Value = '2.989.93'
$ws.Range("E3").Value = '  +3.28%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''563.51'
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").Value = '''138.41'
$ws.Range("E6").Value = '  +11.43%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '''0.520'
$ws.Range("E8").Value = '  +3.69%  '
$ws.Range("D9").Value = '2.986.08'
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("E10").Value = '  +9.29%  '
$ws.Range("D11").Value = '''5.04'
$ws.Range("E11").Value = '  +7.98%  '
$ws.Range("E12").Value = '  +4.84%  '
$ws.Range("D13").Value = '''0.0000230'
$ws.Range("E13").Value = '  +9.67%  '
$ws.Range("D14").Value = '''33.66'
$ws.Range("E14").Value = '  +4.72%  '
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '3.482.22'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").Value = '''7.02'
$ws.Range("E17").Value = '  +7.97%  '
$ws.Range("D18").Value = '2.987.66'
$ws.Range("E18").Value = '  +3.46%  '
$ws.Range("D19").Value = '59.065.12'
$ws.Range("E19").Value = '  +3.33%  '
$ws.Range("D20").Value = '''428.73'
$ws.Range("E20").Value = '  +5.73%  '
$ws.Range("D21").Value = '''13.57'
$ws.Range("E21").Value = '  +6.17%  '
$ws.Range("D22").Value = '''0.717'
$ws.Range("E22").Value = '  +6.87%  '
$ws.Range("D23").Value = '''13.50'
$ws.Range("E23").Value = '  +6.77%  '
$ws.Range("D24").Value = '''7.11'
$ws.Range("E24").Value = '  +4.31%  '
$ws.Range("D25").Value = '''80.52'
$ws.Range("E25").Value = '  +4.02%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").Value = '''2.13'
$ws.Range("E28").Value = '  +11.10%  '
$ws.Range("E29").Value = '  +3.58%  '
$ws.Range("D30").Value = '''7.75'
$ws.Range("E30").Value = '  +7.67%  '
$ws.Range("D31").Value = '''25.66'
$ws.Range("E31").Value = '  +4.12%  '
$ws.Range("D32").Value = '''6.13'
$ws.Range("E32").Value = '  +2.95%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Value = '''0.998'
$ws.Range("E34").Value = '  +8.81%  '
$ws.Range("D35").Value = '0.0₃0768'
$ws.Range("E35").Value = '  +22.58%  '
$ws.Range("D36").Value = '''5.79'
$ws.Range("E36").Value = '  +7.00%  '
$ws.Range("D37").Value = '''2.07'
$ws.Range("E37").Value = '  +3.33%  '
$ws.Range("D38").Value = '''48.98'
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("D39").Value = '''8.67'
$ws.Range("E39").Value = '  +5.11%  '
$ws.Range("D40").Value = '''2.73'
$ws.Range("E40").Value = '  +13.42%  '
$ws.Range("D41").Value = '''397.65'
$ws.Range("E41").Value = '  +10.23%  '
$ws.Range("D42").Value = '''0.0350'
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").Value = '2.744.59'
$ws.Range("E43").Value = '  +4.98%  '
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("D45").Value = '''0.251'
$ws.Range("E45").Value = '  +10.11%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '''123.37'
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("D48").Value = '''0.110'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("D50").Value = '''23.39'
$ws.Range("E50").Value = '  +2.79%  '
$ws.Range("D51").Value = '''32.17'
$ws.Range("E51").Value = '  +19.10%  '
